# The deck's slide master ("theme1.xml") was carrying the "Integral" theme
# colour scheme; the commit swaps it for the stock "Office Theme" palette
# (the second theme part in the package, used by the notes master, already
# holds the "Office" colours). Re-apply the Office Theme's 12 theme colours
# to the presentation's theme via the SlideMaster's ThemeColorScheme, in the
# standard VBA slot order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.

$p = $ppt.ActivePresentation

function HexToRgbLong($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeThemeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$tcs = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $tcs.Colors($i).RGB = HexToRgbLong($officeThemeColors[$i - 1])
}
